# Update the build timestamp embedded in the workbook's version strings.
# Old build: "January 30 2026 16.19.47 EST"
# New build: "February 02 2026 12.49.33 EST"
#
# The string "mines - January 30 (built on January 30 2026 16.19.47 EST)"
# (and the longer citation sentence that embeds it) appears on the "About"
# sheet (A2, A6) and is repeated in the "build_version" column (S) for
# every data row of the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$oldText = "January 30 2026 16.19.47 EST"
$newText = "February 02 2026 12.49.33 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A1:Z24").Replace($oldText, $newText) | Out-Null

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
$wsData.Range("A1:T147").Replace($oldText, $newText) | Out-Null
